$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 and Row 35 had their data swapped (re-sorted source export).
# Only touch cells whose value actually changes to avoid unwanted type coercion
# (e.g. Excel auto-parsing unchanged date-like text cells into date serials).

# --- Row 34 <= original Row 35 values ---
$ws.Range("A34").Value = 131106659
$ws.Range("B34").Value = 57884
$ws.Range("E34").Value = 100109
$ws.Range("F34").Value = 'Tretåig hackspett'
$ws.Range("G34").Value = 'Picoides tridactylus'
$ws.Range("H34").Value = '(Linnaeus, 1758)'
$ws.Range("Q34").Value = 601218
$ws.Range("R34").Value = 6959810
$ws.Range("X34").Value = '2025_0517'
$ws.Range("Z34").Value = '08:33'
$ws.Range("AB34").Value = '08:33'
$ws.Range("AC34").Value = 'På gran i barrblandskog. Ev. liten hackspett, men är med stor sannolikhet tretåig hackspett enligt diskussion med Anders Forsberg.'

# --- Row 35 <= original Row 34 values ---
$ws.Range("A35").Value = 131106648
$ws.Range("B35").Value = 79862
$ws.Range("E35").Value = 6453
$ws.Range("F35").Value = 'Vedskivlav'
$ws.Range("G35").Value = 'Hertelidea botryosa'
$ws.Range("H35").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q35").Value = 601198
$ws.Range("R35").Value = 6959731
$ws.Range("X35").Value = '2025_0528'
$ws.Range("Z35").Value = '09:15'
$ws.Range("AB35").Value = '09:15'
$ws.Range("AC35").Value = 'Brandpåverkad tallstubbe'

# --- Row 38 <= original Row 39 values ---
$ws.Range("A38").Value = 131106646
$ws.Range("B38").Value = 79243
$ws.Range("E38").Value = 6425
$ws.Range("F38").Value = 'Garnlav'
$ws.Range("G38").Value = 'Alectoria sarmentosa'
$ws.Range("H38").Value = '(Ach.) Ach.'
$ws.Range("I38").Value = ''
$ws.Range("J38").ClearContents()
$ws.Range("Q38").Value = 601173
$ws.Range("R38").Value = 6959739
$ws.Range("X38").Value = '2025_0530'
$ws.Range("Z38").Value = '09:36'
$ws.Range("AB38").Value = '09:36'
$ws.Range("AC38").Value = 'tall'
$ws.Range("AX38").Value = 'Måns Svensson'

# --- Row 39 <= original Row 38 values ---
$ws.Range("A39").Value = 131106643
$ws.Range("B39").Value = 78646
$ws.Range("E39").Value = 6437
$ws.Range("F39").Value = 'Blanksvart spiklav'
$ws.Range("G39").Value = 'Calicium denigratum'
$ws.Range("H39").Value = '(Vain.) Tibell'
$ws.Range("I39").Value = "'1"
$ws.Range("J39").Value = 'cm²'
$ws.Range("Q39").Value = 601129
$ws.Range("R39").Value = 6959679
$ws.Range("X39").Value = '2025_0533'
$ws.Range("Z39").Value = '09:55'
$ws.Range("AB39").Value = '09:55'
$ws.Range("AC39").ClearContents()
$ws.Range("AX39").Value = 'Samuel Koont'
